$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1116
$ws.Range("I19").Value = 800
$ws.Range("J19").Value = 1179.2
$ws.Range("K19").Value = 800
$ws.Range("L19").Value = 1179.2
$ws.Range("M19").Value = -625
$ws.Range("N19").Value = -1529.2
$ws.Range("H32").Value = 166667580
$ws.Range("J32").Value = 1131
$ws.Range("L32").Value = 1131
$ws.Range("N32").Value = -1783
$ws.Range("H40").Value = 2155.913
$ws.Range("I40").Value = 1966.0667
$ws.Range("J40").Value = 2511.875
$ws.Range("K40").Value = 1966.0667
$ws.Range("L40").Value = 2511.875
$ws.Range("M40").Value = -1791.0667
$ws.Range("N40").Value = -2861.875
$ws.Range("H43").Value = 2145.484
$ws.Range("I43").Value = 2421.739
$ws.Range("K43").Value = 2421.739
$ws.Range("M43").Value = -2352.739
$ws.Range("H51").Value = 12519.546
$ws.Range("I51").Value = 38000.332
$ws.Range("J51").Value = 2964.25
$ws.Range("K51").Value = 38000.332
$ws.Range("L51").Value = 2964.25
$ws.Range("M51").Value = -37516.332
$ws.Range("N51").Value = -3932.25
$ws.Range("H55").Value = 6114.48
$ws.Range("J55").Value = 7294.8423
$ws.Range("L55").Value = 7294.8423
$ws.Range("N55").Value = -7722.8423
$ws.Range("H113").Value = 201942.2
$ws.Range("I113").Value = 335235
$ws.Range("J113").Value = 2003
$ws.Range("K113").Value = 335235
$ws.Range("L113").Value = 2003
$ws.Range("M113").Value = -331981
$ws.Range("N113").Value = -8511
$ws.Range("H116").Value = 2625.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 7200.1
$ws.Range("I6").Value = 5400.4
$ws.Range("J6").Value = 8999.799999999999
$ws.Range("K6").Value = 5400.4
$ws.Range("L6").Value = 8999.799999999999
$ws.Range("M6").Value = -5227.4
$ws.Range("N6").Value = -9345.799999999999
$ws.Range("H9").Value = 10999
$ws.Range("J9").Value = 10999
$ws.Range("L9").Value = 10999
$ws.Range("N9").Value = -11339
$ws.Range("H20").Value = 10999
$ws.Range("J20").Value = 10999
$ws.Range("L20").Value = 10999
$ws.Range("N20").Value = -11539

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 4490
$ws.Range("I7").Value = 4490
$ws.Range("K7").Value = 4490
$ws.Range("M7").Value = -4377

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 800
$ws.Range("I2").Value = 800
$ws.Range("K2").Value = 800
$ws.Range("M2").Value = -687
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H99").Value = 8821.9375
$ws.Range("I99").Value = 3293.3333
$ws.Range("J99").Value = 10097.77
$ws.Range("K99").Value = 3293.3333
$ws.Range("L99").Value = 10097.77
$ws.Range("M99").Value = -1795.3333
$ws.Range("N99").Value = -13093.77
$ws.Range("H107").Value = 1134.875
$ws.Range("I107").Value = 1575
$ws.Range("J107").Value = 694.75
$ws.Range("K107").Value = 1575
$ws.Range("L107").Value = 694.75
$ws.Range("M107").Value = 345
$ws.Range("N107").Value = -4534.75
$ws.Range("H126").Value = 8821.9375
$ws.Range("I126").Value = 3293.3333
$ws.Range("J126").Value = 10097.77
$ws.Range("K126").Value = 9879.999899999999
$ws.Range("L126").Value = 30293.31
$ws.Range("M126").Value = -7409.999899999999
$ws.Range("N126").Value = -35233.31

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 100000660
$ws.Range("J4").Value = 250001500
$ws.Range("L4").Value = 750004500
$ws.Range("N4").Value = -750004724
$ws.Range("H25").Value = 71429680
$ws.Range("J25").Value = 76924150
$ws.Range("L25").Value = 230772450
$ws.Range("N25").Value = -230772788
$ws.Range("H30").Value = 71429680
$ws.Range("J30").Value = 76924150
$ws.Range("L30").Value = 230772450
$ws.Range("N30").Value = -230772654
$ws.Range("H57").Value = 2674.75
$ws.Range("I57").Value = 849.5
$ws.Range("K57").Value = 2548.5
$ws.Range("M57").Value = -1989.5
$ws.Range("H113").Value = 615.5263
$ws.Range("I113").Value = 598.8570999999999
$ws.Range("J113").Value = 625.25
$ws.Range("K113").Value = 1796.5713
$ws.Range("L113").Value = 1875.75
$ws.Range("M113").Value = 373.4287000000002
$ws.Range("N113").Value = -6215.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 508380.9
$ws.Range("I5").Value = 5000000
$ws.Range("J5").Value = 9312.111000000001
$ws.Range("K5").Value = 5000000
$ws.Range("L5").Value = 9312.111000000001
$ws.Range("M5").Value = -4999888
$ws.Range("N5").Value = -9536.111000000001
$ws.Range("H12").Value = 5349396
$ws.Range("J12").Value = 4420.8
$ws.Range("L12").Value = 4420.8
$ws.Range("N12").Value = -4700.8
$ws.Range("H132").Value = 3300.9473
$ws.Range("I132").Value = 2759.5
$ws.Range("J132").Value = 3902.5557
$ws.Range("K132").Value = 8278.5
$ws.Range("L132").Value = 11707.6671
$ws.Range("M132").Value = -5748.5
$ws.Range("N132").Value = -16767.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 800
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -1590
$ws.Range("H26").Value = 10000
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H27").Value = 800
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -1214
$ws.Range("H132").Value = 4471.5
$ws.Range("I132").Value = 4826.4614
$ws.Range("J132").Value = 3548.6
$ws.Range("K132").Value = 14479.3842
$ws.Range("L132").Value = 10645.8
$ws.Range("M132").Value = -11949.3842
$ws.Range("N132").Value = -15705.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 168117.83
$ws.Range("I81").Value = 200511
$ws.Range("J81").Value = 144979.86
$ws.Range("K81").Value = 401022
$ws.Range("L81").Value = 289959.72
$ws.Range("M81").Value = -399961
$ws.Range("N81").Value = -292081.72
$ws.Range("H84").Value = 168117.83
$ws.Range("I84").Value = 200511
$ws.Range("J84").Value = 144979.86
$ws.Range("K84").Value = 2005110
$ws.Range("L84").Value = 1449798.6
$ws.Range("M84").Value = -1999806
$ws.Range("N84").Value = -1460406.6
